$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58, shifting rows 58:164 down to 59:165
$ws.Rows(58).Insert()

# Populate the newly inserted row 58 with the new weekly price record
$ws.Range("A58").Value = 7
$ws.Range("B58").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C58").Value = "Ñuble"
$ws.Range("D58").Value = 44469
$ws.Range("E58").Value = 16
$ws.Range("F58").Value = 100112008
$ws.Range("G58").Value = "Coliflor"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 300
$ws.Range("K58").Value = 700
$ws.Range("L58").Value = 750
$ws.Range("M58").Value = 725
$ws.Range("N58").Value = "$/unidad"
$ws.Range("O58").Value = "Región del Maule"
$ws.Range("P58").Value = 725
$ws.Range("Q58").Value = 1
$ws.Range("R58").Value = "Hortaliza"
